$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new run-log row (row 58) mirroring the existing entries.
$row = 58

# Match the formatting (style) of the preceding data row first.
$ws.Range("A57:H57").Copy()
$ws.Range("A58:H58").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = "2025-08-26 03:52:06 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-26 09:22:06 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""
